$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.778.16'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.303.73'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '156.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +15,516.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.55'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '96.89'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +5.10%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.500'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '35.71'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +8.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0811'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.62%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.79'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.660.30'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.63'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.296.27'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.801'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.625.13'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.86'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +5.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0924'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.44'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '245.07'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.63'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.49%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.38'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.25'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +8.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.74'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.67'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.11%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0759'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.41%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.45'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.72%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.47%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.40%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.23'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.13'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.026.63'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +10.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0286'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.30'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.06'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.37%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.62%  '
